$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 655; Excel shifts rows 655-737 down to 656-738
# and keeps formatting from the row above (matches the diff's dimension
# change from A1:T737 to A1:T738 and every subsequent row's content moving
# down by one).
$ws.Rows.Item(655).Insert()

# Columns A, B, C, E, F, G, H, I, J are constant across this block
# (Macroferia Regional de Talca / Maule / Fruta / Citricos / Naranja), so
# copy them down from the row above rather than retyping literals.
$ws.Range("A655").Value = $ws.Range("A656").Value2
$ws.Range("B655").Value = $ws.Range("B656").Value2
$ws.Range("C655").Value = $ws.Range("C656").Value2
$ws.Range("E655").Value = $ws.Range("E656").Value2
$ws.Range("F655").Value = $ws.Range("F656").Value2
$ws.Range("G655").Value = $ws.Range("G656").Value2
$ws.Range("H655").Value = $ws.Range("H656").Value2
$ws.Range("I655").Value = $ws.Range("I656").Value2
$ws.Range("J655").Value = $ws.Range("J656").Value2

# New data for the inserted row.
$ws.Range("D655").Value = 44918
$ws.Range("K655").Value = "Navel Late"
$ws.Range("L655").Value = "Primera"
$ws.Range("M655").Value = 320
$ws.Range("N655").Value = 10000
$ws.Range("O655").Value = 10000
$ws.Range("P655").Value = 10000
$ws.Range("Q655").Value = '$/bandeja 15 kilos granel'
$ws.Range("R655").Value = "Provincia de Melipilla"
$ws.Range("S655").Value = 667
$ws.Range("T655").Value = 15
